$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$entryDate = Get-Date -Year 2019 -Month 12 -Day 7 -Hour 0 -Minute 0 -Second 0

# Row 2
$ws.Range("A2").NumberFormat = "YYYY-MM-DD"
$ws.Range("A2").Value = $entryDate
$ws.Range("B2").Style = "Normal"
$ws.Range("B2").Value = "https://www.google.com/search?ei=fRZjXOTcNoOYafn9rKgG"

# Row 3
$ws.Range("A3").NumberFormat = "YYYY-MM-DD"
$ws.Range("A3").Value = $entryDate
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").Value = "https://www.youtube.com/watch?v=zA0eqkqwaKE"

# Row 4
$ws.Range("A4").NumberFormat = "YYYY-MM-DD"
$ws.Range("A4").Value = $entryDate
$ws.Range("B4").Style = "Normal"
$ws.Range("B4").Value = "https://olympus.greatlearning.in/login"
$ws.Range("C4").Style = "Normal"
$ws.Range("C4").Value = "http://127.0.0.1:5000/cV63QT"

# Row 5
$ws.Range("A5").NumberFormat = "YYYY-MM-DD"
$ws.Range("A5").Value = $entryDate
$ws.Range("B5").Style = "Normal"
$ws.Range("B5").Value = "https://zoom.us/wc/leave?meetingNumber=442013660"
$ws.Range("C5").Style = "Normal"
$ws.Range("C5").Value = "http://127.0.0.1:5000/shxMzb"
